$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''28.925.30'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  -2.49%  '
$ws.Range("D3").Value = '''1.895.49'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  -5.41%  '
$ws.Range("E4").Value = '  -0.39%  '
$ws.Range("D5").Value = '''323.73'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -1.50%  '
$ws.Range("E6").Value = '  -0.47%  '
$ws.Range("D7").Value = '''0.4601'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -2.14%  '
$ws.Range("D8").Value = '''0.3821'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -3.60%  '
$ws.Range("E9").Value = '  -3.03%  '
$ws.Range("D10").Value = '''0.07733'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -2.95%  '
$ws.Range("D11").Value = '''0.9712'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -3.54%  '
$ws.Range("D12").Value = '''22.17'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -3.08%  '
$ws.Range("D13").Value = '''1.892.62'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -5.53%  '
$ws.Range("D14").Value = '''5.699'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -3.42%  '
$ws.Range("D15").Value = '''6.968'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -4.22%  '
$ws.Range("D16").Value = '''0.07036'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -1.84%  '
$ws.Range("E17").Value = '  -0.47%  '
$ws.Range("D18").Value = '''83.60'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -6.33%  '
$ws.Range("D19").Value = '''0.000009537'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -4.76%  '
$ws.Range("D20").Value = '''16.74'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -3.93%  '
$ws.Range("E21").Value = '  -0.50%  '
$ws.Range("D22").Value = '''28.926.39'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -2.71%  '
$ws.Range("D23").Value = '''5.333'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -3.97%  '
$ws.Range("D25").Value = '''2.133.97'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -5.08%  '
$ws.Range("E26").Value = '  -2.81%  '
$ws.Range("D27").Value = '''156.19'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -1.67%  '
$ws.Range("D28").Value = '''19.13'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -3.24%  '
$ws.Range("D29").Value = '''5.620'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -5.95%  '
$ws.Range("D30").Value = '''117.51'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -2.73%  '
$ws.Range("D31").Value = '''1.817'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -7.98%  '
$ws.Range("D32").Value = '''0.09268'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -2.31%  '
$ws.Range("D33").Value = '''0.8548'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -5.18%  '
$ws.Range("D34").Value = '''5.094'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -4.00%  '
$ws.Range("D35").Value = '''1.242'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -7.87%  '
$ws.Range("D36").Value = '''3.012'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -5.78%  '
$ws.Range("D37").Value = '''0.05704'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -2.53%  '
$ws.Range("D38").Value = '''1.147'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -2.55%  '
$ws.Range("E39").Value = '  -0.44%  '
$ws.Range("D40").Value = '''0.02040'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -4.70%  '
$ws.Range("B41").Value = 'TheSandbox'
$ws.Range("C41").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D41").Value = '''0.5512'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -4.68%  '
$ws.Range("B42").Value = 'FraxShare'
$ws.Range("C42").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D42").Value = '''7.436'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -6.18%  '
$ws.Range("D43").Value = '''0.1750'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -4.20%  '
$ws.Range("D44").Value = '''0.000002882'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -16.02%  '
$ws.Range("D45").Value = '''9.269'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -6.32%  '
$ws.Range("D46").Value = '''2.702'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +1.49%  '
$ws.Range("D47").Value = '''0.5193'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -3.90%  '
$ws.Range("D48").Value = '''11.34'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -6.85%  '
$ws.Range("D49").Value = '''2.083'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -4.62%  '
$ws.Range("D50").Value = '''0.06818'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -2.58%  '
$ws.Range("D51").Value = '''111.62'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -2.80%  '
